$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.193.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.587.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +1.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.18%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.814.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.582.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.234.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.396.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.63%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("E36").Value = "  -7.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("E39").Value = "  +8.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.540"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.981"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.724.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("E51").Value = "  -0.64%  "
